$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Type ID / Type Name values between B2 and B3
$ws.Range("B2").Value = "Z.Sandbox.HWDBUnitTest.biff"
$ws.Range("B3").Value = "Z00100300022"

# Update generated row (row 9) with new values
$ws.Range("B9").Value = "03532066-3DF8-44E2-9E77-385DEB7FFAE0"
$ws.Range("C9").Value = "generated 2023-12-12 00:17:40"
$ws.Range("D9").Value = 116.89
$ws.Range("E9").Value = 97.73999999999999
$ws.Range("F9").Value = 110.21
$ws.Range("G9").Value = "3539CBF4-04D5-4BDD-AAE3-61CBB5F3A2CA"
$ws.Range("H9").Value = "B93E8AA7-F6FE-4384-A722-6A4359B6A0A4"
